$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "AddCustomerTest"

# Header row (first three columns)
$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "lastname"
$ws.Range("C1").Value = "postcode"

# Data row (first three columns)
$ws.Range("A2").Value = "Harry"
$ws.Range("B2").Value = "Potter"
$ws.Range("C2").Value = "HG314"

# Fourth column added afterwards
$ws.Range("D1").Value = "alerttext"
$ws.Range("D2").Value = "Customer added successfully"

# Set selection to D2 as shown in the diff
$ws.Range("D2").Select()
